$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.863.11"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "3.481.86"
$ws.Range("E3").Value = "  -0.01%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.95"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.64"
$ws.Range("E6").Value = "  -2.32%  "

$ws.Range("D7").Value = "3.479.25"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.93"
$ws.Range("E11").Value = "  +5.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.416"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").Value = "4.068.74"
$ws.Range("E14").Value = "  -0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.01"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "3.472.60"
$ws.Range("E16").Value = "  -0.77%  "

$ws.Range("D17").Value = "66.932.02"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("E18").Value = "  -0.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  +6.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.29"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.32"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("E23").Value = "  -3.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.21"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "3.620.47"
$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  -3.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("E29").Value = "  -4.23%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  -4.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.164"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.32"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("E35").Value = "  -2.83%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("E38").Value = "  -7.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.43"
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("E41").Value = "  -0.41%  "

$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.96"
$ws.Range("E44").Value = "  -12.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.18"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.66"
$ws.Range("E46").Value = "  -7.92%  "

$ws.Range("E47").Value = "  -4.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.29"
$ws.Range("E48").Value = "  -3.06%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  -3.32%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.972"
$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.242"
$ws.Range("E51").Value = "  -1.63%  "
